# Add the "Michele Leonardi " team as a new row (row 38) at the bottom of Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A38").Value = "Michele Leonardi "
$ws.Range("B38").Value = "Elia Battisti | U.SGUARNA"
$ws.Range("C38").Value = "Filippo Benetti | I Magnifici"
$ws.Range("D38").Value = "Matteo Diener | U.SGUARNA"
$ws.Range("E38").Value = "Luca Tonolli | Rita Levi’s"
$ws.Range("F38").Value = "Sayf Brik | A.C.DENTI"
